$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 6220
$ws.Range("L3").Value = 6750
$ws.Range("L4").Value = 1669
$ws.Range("L5").Value = 399
$ws.Range("L6").Value = 5551
$ws.Range("L7").Value = 20589

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L2").Value = 182
$ws.Range("L8").Value = 1360
$ws.Range("L15").Value = 168
$ws.Range("L16").Value = 47
$ws.Range("L17").Value = 36
$ws.Range("L19").Value = 556
$ws.Range("L20").Value = 526
$ws.Range("L22").Value = 67
$ws.Range("L23").Value = 218
$ws.Range("L24").Value = 60
$ws.Range("L25").Value = 125
$ws.Range("L27").Value = 177
$ws.Range("L29").Value = 1150
$ws.Range("L35").Value = 26
$ws.Range("L37").Value = 792
$ws.Range("L42").Value = 657
$ws.Range("L43").Value = 153
$ws.Range("L47").Value = 146
$ws.Range("L48").Value = 270
$ws.Range("L50").Value = 100
$ws.Range("L53").Value = 230
$ws.Range("L55").Value = 219
$ws.Range("L63").Value = 60
$ws.Range("L66").Value = 58
$ws.Range("L67").Value = 715
$ws.Range("L72").Value = 83
$ws.Range("L73").Value = 162
$ws.Range("L76").Value = 319
$ws.Range("L79").Value = 568
$ws.Range("L85").Value = 1024
$ws.Range("L87").Value = 56
$ws.Range("L90").Value = 217
$ws.Range("L91").Value = 277
$ws.Range("L94").Value = 251
$ws.Range("L95").Value = 290
$ws.Range("L96").Value = 228
$ws.Range("L101").Value = 20589

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L3").Value = 65
$ws.Range("L7").Value = 228

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L2").Value = 309
$ws.Range("L6").Value = 212
$ws.Range("L7").Value = 1024

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("L2").Value = 69
$ws.Range("L3").Value = 58
$ws.Range("L7").Value = 230

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 412
$ws.Range("L3").Value = 483
$ws.Range("L7").Value = 1360

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L2").Value = 104
$ws.Range("L7").Value = 290

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L2").Value = 240
$ws.Range("L3").Value = 281
$ws.Range("L7").Value = 792

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L3").Value = 280
$ws.Range("L6").Value = 165
$ws.Range("L7").Value = 715

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L3").Value = 443
$ws.Range("L6").Value = 281
$ws.Range("L7").Value = 1150

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L6").Value = 110
$ws.Range("L7").Value = 270

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L2").Value = 201
$ws.Range("L4").Value = 26
$ws.Range("L7").Value = 556

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L3").Value = 65
$ws.Range("L7").Value = 319

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L3").Value = 226
$ws.Range("L6").Value = 186
$ws.Range("L7").Value = 657

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("L6").Value = 61
$ws.Range("L7").Value = 219

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("L4").Value = 4
$ws.Range("L7").Value = 60

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("L6").Value = 54
$ws.Range("L7").Value = 218

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("L6").Value = 36
$ws.Range("L7").Value = 277

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L2").Value = 176
$ws.Range("L6").Value = 152
$ws.Range("L7").Value = 568

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L2").Value = 164
$ws.Range("L7").Value = 526

$ws = $wb.Worksheets.Item("Burnside")
$ws.Range("L3").Value = 14
$ws.Range("L7").Value = 36

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L3").Value = 60
$ws.Range("L7").Value = 251

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("L5").Value = 6
$ws.Range("L6").Value = 18
$ws.Range("L7").Value = 125

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("L3").Value = 50
$ws.Range("L7").Value = 146

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("L4").Value = 14
$ws.Range("L7").Value = 168

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("L3").Value = 27
$ws.Range("L7").Value = 100

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("L6").Value = 17
$ws.Range("L7").Value = 58

$ws = $wb.Worksheets.Item("Gold Coast")
$ws.Range("L3").Value = 6
$ws.Range("L7").Value = 26

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("L4").Value = 16
$ws.Range("L7").Value = 162

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("L3").Value = 59
$ws.Range("L6").Value = 47
$ws.Range("L7").Value = 182

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("L3").Value = 51
$ws.Range("L7").Value = 177

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("L2").Value = 72
$ws.Range("L7").Value = 217

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("L6").Value = 46
$ws.Range("L7").Value = 153

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("L3").Value = 26
$ws.Range("L7").Value = 67

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("L4").Value = 13
$ws.Range("L7").Value = 83

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("L3").Value = 14
$ws.Range("L7").Value = 56

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("L6").Value = 31
$ws.Range("L7").Value = 47
